$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the check digit (dígito verificador) of Juan Pérez's RUT,
# which was "12345678-9" and should be "12345678-5".
$ws.Range("C3").Value = "12345678-5"

# Reflect the last clicked/selected cell as seen in the saved file.
$ws.Range("C3").Select()
